# Applies the "Marketing Strategy + PR" update to Paula Hernandez's resume.
#
# Helper: run a plain (format-agnostic) Find & Replace against the whole
# document story, matching exactly once. Throws if the text wasn't found,
# so mistakes are not silent.
function Replace-Text($old, $new) {
    $range = $word.ActiveDocument.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute could not locate: $old"
    }
}

$d = $word.ActiveDocument

# 1) Professional summary - bolded role title: add "Public Relations"
Replace-Text `
    "Marketing Strategy and Strategic Partnerships function at MyBambu" `
    "Marketing Strategy, Public Relations, and Strategic Partnerships function at MyBambu"

# 2) Professional summary - "integrated marketing and communication strategies"
Replace-Text `
    "—I create integrated marketing and communication strategies across influencers, traditional media, brand partnerships, events, and digital content." `
    "—I create integrated marketing, PR, and communication strategies across influencers, traditional media, brand partnerships, events, and digital content."

# 3) Professional summary - "from influencer programs..." lead-in to "I specialize in"
Replace-Text `
    "—from influencer programs and ambassador networks to national media campaigns and community sponsorships. I specialize in" `
    "—from PR and media relations to influencer programs, ambassador networks, national media campaigns, and community sponsorships. I specialize in"

# 4) Professional summary - bolded "strategic storytelling ..." run is shortened
#    (the descriptive tail that used to live in this bold run moves into the
#    following plain-text run, see step 5).
Replace-Text `
    "strategic storytelling that connects with community, strengthens brand trust, and elevates visibility at scale" `
    "strategic storytelling, media relations, crisis communications, and brand positioning"

# 5) Professional summary - plain run right after the bold run above absorbs
#    the former bold tail text (now un-bolded) plus new PR-related content.
Replace-Text `
    ". My experience spans building and executing comprehensive marketing initiatives including national campaigns, multi-city tours, influencer ecosystems, ambassador programs, and cross-functional brand partnerships—generating" `
    " that connects with community, strengthens brand trust, and elevates visibility at scale. My experience spans building and executing comprehensive marketing and PR initiatives including national campaigns, press releases, media placements, multi-city tours, influencer ecosystems, ambassador programs, and cross-functional brand partnerships—generating"

# 6) Professional summary - "Through these integrated marketing efforts..."
Replace-Text `
    ". Through these integrated marketing efforts across all channels, I have contributed to" `
    ". Through these integrated marketing and PR efforts across all channels, I have contributed to"

# 7) Experience section heading line
Replace-Text `
    "Founding Leader — Built Marketing Strategy & Partnerships Department from Scratch" `
    "Founding Leader — Built Marketing Strategy, PR & Partnerships Department from Scratch"

# 8) Bullet label: "Marketing Strategy & Department Development:"
Replace-Text `
    "Marketing Strategy & Department Development:" `
    "Marketing Strategy, PR & Department Development:"

# 9) Bullet body text describing the department build-out
Replace-Text `
    "Built the entire Marketing Strategy and Strategic Partnerships function from the ground up, creating every initiative from scratch. Defined comprehensive strategy, processes, workflows, messaging guidelines, KPIs, and brand communication pillars. Led national marketing and PR strategy targeting Hispanic audiences across the U.S., ensuring alignment with brand mission, cultural values, and growth goals. Transformed a brand with zero marketing infrastructure into a recognized leader in the Hispanic FinTech market." `
    "Built the entire Marketing Strategy, Public Relations, and Strategic Partnerships function from the ground up, creating every initiative from scratch. Defined comprehensive marketing and PR strategy, processes, workflows, messaging guidelines, KPIs, and brand communication pillars. Led national marketing, PR, and media relations strategy targeting Hispanic audiences across the U.S., ensuring alignment with brand mission, cultural values, and growth goals. Transformed a brand with zero marketing and PR infrastructure into a recognized leader in the Hispanic FinTech market."

Write-Output "All replacements applied."
